# Apply equity_curve.xlsx update:
#   Initial equity rebased from 1600 -> 1000 (and the trailing backtest
#   rows recomputed for the new initial cash amount).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-43: constant equity of 1600 -> 1000
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 1).Value = 1000
}

# Row 44: equity changes, drawdown (column B) stays 0
$ws.Cells.Item(44, 1).Value = 1089.4233935803

# Rows 45-53: equity (A) and drawdown pct (B) recomputed
$ws.Cells.Item(45, 1).Value = 970.0966492703003
$ws.Cells.Item(45, 2).Value = 0.1095320194271231

$ws.Cells.Item(46, 1).Value = 939.9616580703004
$ws.Cells.Item(46, 2).Value = 0.1371934331415503

$ws.Cells.Item(47, 1).Value = 952.3076986503002
$ws.Cells.Item(47, 2).Value = 0.1258607954795065

$ws.Cells.Item(48, 1).Value = 971.2272064303003
$ws.Cells.Item(48, 2).Value = 0.1084942620532113

$ws.Cells.Item(49, 1).Value = 962.5134024703002
$ws.Cells.Item(49, 2).Value = 0.1164928088178102

$ws.Cells.Item(50, 1).Value = 977.7303850503004
$ws.Cells.Item(50, 2).Value = 0.1025248853551143

$ws.Cells.Item(51, 1).Value = 1030.0547879703
$ws.Cells.Item(51, 2).Value = 0.05449543855937389

$ws.Cells.Item(52, 1).Value = 1043.820459870301
$ws.Cells.Item(52, 2).Value = 0.0418596975048694

$ws.Cells.Item(53, 1).Value = 1043.820459870301
$ws.Cells.Item(53, 2).Value = 0.0418596975048694

$wb.Save()
